$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 148; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -is [double]) {
        $cell.Value2 = -1 * $val
    }
}
